# Update the "取得日時" (acquisition timestamp) column on the "ランサーズ" sheet
# for all existing data rows from the old run time to the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-12-05 06:28:47"
$newTimestamp = "2025-12-05 06:37:53"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Text -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
